# M5 Client.pptx - "overhauled client code to use hottowel template"
#
# Inserts 5 new slides between the existing "Module Overview"-style slides
# (positions 1-5, untouched) and the existing "Resources"/"Summary"/blank
# slides (which shift from positions 6-8 down to 11-13):
#   6.  Options for Getting Started (bulleted list of options)
#   7.  Yeoman generator-angular
#   8.  Cloning angular/angular-seed
#   9.  Starting with the ng.NET project template
#   10. Starting with HotTowel by John Papa via Nuget

$p = $ppt.ActivePresentation

# --- Slide 6: "Options for Getting Started" (Title + Content placeholder) ---
$s6 = $p.Slides.Add(6, 2)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Options for Getting Started"
$s6.Shapes.Item(2).TextFrame.TextRange.Text = "angular from scratch`ryeoman generator-angular`rangular/angular-seed`rng.NET (VS project template)`rHot Towel by John Papa (via Nuget)"

# --- Slide 7: "Yeoman generator-angular" (Title only) ---
$s7 = $p.Slides.Add(7, 11)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Yeoman generator-angular"

# --- Slide 8: "Cloning angular/angular-seed" (Title only) ---
$s8 = $p.Slides.Add(8, 11)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "Cloning angular/angular-seed"

# --- Slide 9: "Starting with the ng.NET project template" (Title only) ---
$s9 = $p.Slides.Add(9, 11)
$s9.Shapes.Item(1).TextFrame.TextRange.Text = "Starting with the ng.NET project template"

# --- Slide 10: "Starting with HotTowel by John Papa via Nuget" (Title only) ---
$s10 = $p.Slides.Add(10, 11)
$s10.Shapes.Item(1).TextFrame.TextRange.Text = "Starting with HotTowel by John Papa via Nuget"
